$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'26.983.25"
$ws.Range("E2").Value = "  -2.04%  "

$ws.Range("D3").Value = "'1.821.76"
$ws.Range("E3").Value = "  -0.93%  "

$ws.Range("E4").Value = "  -0.61%  "

$ws.Range("D5").Value = "'310.58"
$ws.Range("E5").Value = "  -1.97%  "

$ws.Range("E6").Value = "  -0.68%  "

$ws.Range("D7").Value = "'0.4251"
$ws.Range("E7").Value = "  -1.25%  "

$ws.Range("D8").Value = "'0.3655"
$ws.Range("E8").Value = "  -1.60%  "

$ws.Range("D9").Value = "'0.07223"
$ws.Range("E9").Value = "  -0.97%  "

$ws.Range("D10").Value = "'0.8408"
$ws.Range("E10").Value = "  -3.34%  "

$ws.Range("D11").Value = "'20.54"
$ws.Range("E11").Value = "  -3.26%  "

$ws.Range("D12").Value = "'1.842.84"
$ws.Range("E12").Value = "  +0.16%  "

$ws.Range("D13").Value = "'6.639"
$ws.Range("E13").Value = "  -1.10%  "

$ws.Range("D14").Value = "'0.07062"
$ws.Range("E14").Value = "  -0.80%  "

$ws.Range("E15").Value = "  -1.99%  "

$ws.Range("D16").Value = "'89.28"
$ws.Range("E16").Value = "  +0.75%  "

$ws.Range("D17").Value = "'1.001"
$ws.Range("E17").Value = "  -0.82%  "

$ws.Range("D18").Value = "'0.000008738"
$ws.Range("E18").Value = "  -2.41%  "

$ws.Range("E19").Value = "  -0.59%  "

$ws.Range("D20").Value = "'14.84"
$ws.Range("E20").Value = "  -3.15%  "

$ws.Range("D21").Value = "'27.035.86"
$ws.Range("E21").Value = "  -1.88%  "

$ws.Range("D22").Value = "'5.117"
$ws.Range("E22").Value = "  -1.27%  "

$ws.Range("E23").Value = "  -1.59%  "

$ws.Range("D24").Value = "'2.047.00"
$ws.Range("E24").Value = "  -1.03%  "

$ws.Range("D25").Value = "'1.976"
$ws.Range("E25").Value = "  -1.71%  "

$ws.Range("D26").Value = "'150.95"
$ws.Range("E26").Value = "  -2.28%  "

$ws.Range("D27").Value = "'2.216"
$ws.Range("E27").Value = "  +2.47%  "

$ws.Range("E28").Value = "  -2.05%  "

$ws.Range("D29").Value = "'5.208"
$ws.Range("E29").Value = "  -2.00%  "

$ws.Range("D30").Value = "'116.60"
$ws.Range("E30").Value = "  -0.79%  "

$ws.Range("D31").Value = "'0.08714"
$ws.Range("E31").Value = "  -1.94%  "

$ws.Range("D32").Value = "'1.173"
$ws.Range("E32").Value = "  -3.24%  "

$ws.Range("D33").Value = "'0.7334"
$ws.Range("E33").Value = "  -4.85%  "

$ws.Range("D34").Value = "'2.901"
$ws.Range("E34").Value = "  -0.31%  "

$ws.Range("D35").Value = "'4.410"
$ws.Range("E35").Value = "  -2.14%  "

$ws.Range("D36").Value = "'1.000"
$ws.Range("E36").Value = "  -0.77%  "

$ws.Range("E37").Value = "  -3.46%  "

$ws.Range("D38").Value = "'0.01937"
$ws.Range("E38").Value = "  -1.38%  "

$ws.Range("D39").Value = "'0.05205"
$ws.Range("E39").Value = "  -1.85%  "

$ws.Range("D40").Value = "'7.203"
$ws.Range("E40").Value = "  +0.49%  "

$ws.Range("D41").Value = "'2.857"
$ws.Range("E41").Value = "  -0.87%  "

$ws.Range("D42").Value = "'0.1681"
$ws.Range("E42").Value = "  +0.11%  "

$ws.Range("D43").Value = "'0.5097"
$ws.Range("E43").Value = "  -0.16%  "

$ws.Range("D44").Value = "'8.515"
$ws.Range("E44").Value = "  -2.51%  "

$ws.Range("D45").Value = "'10.54"
$ws.Range("E45").Value = "  -0.47%  "

$ws.Range("D46").Value = "'1.951"
$ws.Range("E46").Value = "  +6.08%  "

$ws.Range("D47").Value = "'0.4721"
$ws.Range("E47").Value = "  -0.20%  "

$ws.Range("D48").Value = "'105.64"
$ws.Range("E48").Value = "  -1.04%  "

$ws.Range("D49").Value = "'0.9998"
$ws.Range("E49").Value = "  -0.81%  "

$ws.Range("D50").Value = "'0.06310"
$ws.Range("E50").Value = "  -2.04%  "

$ws.Range("E51").Value = "  -1.94%  "
